# updates per Eric K. comments
#
# 1) Refresh the auto date field (datetimeFigureOut) on the slide master
#    and every slide layout from 6/9/2021 -> 7/1/2021.
# 2) "Laboratory environment" textbox -> "Laboratory and curatorial environment"
#    (middle phrase gets Slack-paste-style run formatting) + reposition/resize.
# 3) The connector glued to that textbox gets a shorter/shallower arrow.
# 4) "Historic human occupation site" -> "Site of past human activities"
#    (again with the Slack-paste-style run formatting on the bulk of the text)
#    + reposition/resize.
# 5) "e.g. lab blank" textbox nudges position.
# 6) The connector glued to the Laboratory textbox gets a shorter arrow.
# 7) "A place where humans have been..." caption grows and gains a
#    parenthetical about prehistoric hominids.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholders: slide master + all custom layouts.
# ---------------------------------------------------------------------------
function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDate = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if ($isDate) {
            $shp.TextFrame.TextRange.Text = "7/1/2021"
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes
for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $lay = $master.CustomLayouts.Item($j)
    Update-DatePlaceholders $lay.Shapes
}

# ---------------------------------------------------------------------------
# Working slide.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# 2) "Laboratory environment" -> "Laboratory and curatorial environment"
# ---------------------------------------------------------------------------
$lab = $s.Shapes.Item("TextBox 125")
$labTr = $lab.TextFrame.TextRange
$labTr.Text = "Laboratory and curatorial environment"

$labMid = $labTr.Characters(12, 14)  # "and curatorial"
$labMid.Font.Italic = 0
$labMid.Font.Color.RGB = 0x1D1C1D
$labMid.Font.Name = "Slack-Lato"
$labMid.Font.Shadow = 0

$lab.Left = 701.7119750976562
$lab.Top = 373.0556335449219
$lab.Height = 72.70311737060547

# ---------------------------------------------------------------------------
# 3) Connector glued to the Laboratory textbox (endCxn id=126).
# ---------------------------------------------------------------------------
$connLab = $s.Shapes.Item("Straight Arrow Connector 145")
$connLab.Width = 54.4102783203125

# ---------------------------------------------------------------------------
# 4) "Historic human occupation site" -> "Site of past human activities"
# ---------------------------------------------------------------------------
$hist = $s.Shapes.Item("TextBox 336")
$histTr = $hist.TextFrame.TextRange
$histTr.Text = "Site of past human activities"

$histRest = $histTr.Characters(2, 29)  # "ite of past human activities"
$histRest.Font.Color.RGB = 0x1D1C1D
$histRest.Font.Name = "Slack-Lato"
$histRest.Font.Italic = 0
$histRest.Font.Shadow = 0

$histFirst = $histTr.Characters(1, 1)  # "S"
$histFirst.Font.Color.RGB = 0x1D1C1D
$histFirst.Font.Name = "Slack-Lato"

$hist.Left = 915.9688720703125
$hist.Width = 141.9190216064453

# ---------------------------------------------------------------------------
# 5) "e.g. lab blank" textbox repositions.
# ---------------------------------------------------------------------------
$labBlank = $s.Shapes.Item("TextBox 162")
$labBlank.Left = 719.7681274414062
$labBlank.Top = 442.98516845703125

# ---------------------------------------------------------------------------
# 6) Connector glued to the Laboratory textbox's start side (stCxn id=223,
#    endCxn id=126) shrinks.
# ---------------------------------------------------------------------------
$connLab2 = $s.Shapes.Item("Straight Arrow Connector 225")
$connLab2.Width = 59.37681198120117
$connLab2.Height = 36.1298828125

# ---------------------------------------------------------------------------
# 7) "A place where humans have been..." caption text + height.
# ---------------------------------------------------------------------------
$humans = $s.Shapes.Item("TextBox 128")
$humans.TextFrame.TextRange.Text = "A place where humans (including related prehistoric hominids)  have been and left evidence of their activity"
$humans.Height = 55.739017486572266

Write-Host "edit complete"
